# Rename the three logo pictures (header "first page" BTec logo + both
# footers' Pearson logo) so each InlineShape's Name no longer collides
# with its sibling logo's filename:
#   BTec logo (header)    : image2.jpg -> image1.jpg
#   Pearson logo (footers): image1.png -> image2.png
#
# NOTE: InlineShape.Name is a write-only-ish property in this runtime (its
# getter does not reflect the name already baked into the document XML),
# so shapes are identified via AlternativeText (which *does* read back
# correctly) rather than by their current Name.

$d = $word.ActiveDocument

$btecAlt    = "BTec_Logo-Orange"
$pearsonAlt = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"

foreach ($sec in $d.Sections) {

    # --- Headers ---------------------------------------------------
    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $hdr = $sec.Headers.Item($i)
        if (-not $hdr.Exists) { continue }
        for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
            $shp = $hdr.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -eq $btecAlt) {
                $shp.Name = "image1.jpg"
            }
            elseif ($shp.AlternativeText -eq $pearsonAlt) {
                $shp.Name = "image2.png"
            }
        }
    }

    # --- Footers -----------------------------------------------------
    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $ftr = $sec.Footers.Item($i)
        if (-not $ftr.Exists) { continue }
        for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
            $shp = $ftr.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -eq $btecAlt) {
                $shp.Name = "image1.jpg"
            }
            elseif ($shp.AlternativeText -eq $pearsonAlt) {
                $shp.Name = "image2.png"
            }
        }
    }
}
